$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price (D) and volume (E) cells keep their original text formatting
# (e.g. "1.00", "97.540.87") instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "97.540.87"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "3.700.99"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "2.17"
$ws.Range("E5").Value = "  +12.91%  "
$ws.Range("D6").Value = "237.10"
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("D7").Value = "655.61"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").Value = "0.446"
$ws.Range("E8").Value = "  +3.96%  "
$ws.Range("D9").Value = "1.13"
$ws.Range("E9").Value = "  +3.22%  "
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "3.697.44"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Value = "0.0000313"
$ws.Range("E12").Value = "  +15.26%  "
$ws.Range("D13").Value = "44.82"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").Value = "0.207"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "4.392.52"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "96.976.31"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "8.87"
$ws.Range("E18").Value = "  -2.42%  "
$ws.Range("D19").Value = "3.698.80"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").Value = "13.02"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "18.75"
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("D22").Value = "0.539"
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("D23").Value = "523.40"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").Value = "3.45"
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("D25").Value = "0.0000222"
$ws.Range("E25").Value = "  +7.93%  "
$ws.Range("D26").Value = "117.78"
$ws.Range("E26").Value = "  +14.43%  "
$ws.Range("D27").Value = "6.91"
$ws.Range("E27").Value = "  -2.70%  "
$ws.Range("D28").Value = "0.211"
$ws.Range("E28").Value = "  +24.65%  "
$ws.Range("D29").Value = "13.46"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "12.70"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "0.189"
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("D35").Value = "32.92"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "0.595"
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("D38").Value = "637.08"
$ws.Range("E38").Value = "  -3.42%  "
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("D42").Value = "0.499"
$ws.Range("E42").Value = "  +12.28%  "
$ws.Range("D43").Value = "6.82"
$ws.Range("E43").Value = "  -4.99%  "
$ws.Range("D44").Value = "40.10"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("D45").Value = "2.00"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").Value = "0.959"
$ws.Range("E46").Value = "  -1.98%  "
$ws.Range("D47").Value = "0.0453"
$ws.Range("E47").Value = "  -1.82%  "
$ws.Range("D48").Value = "2.37"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").Value = "8.80"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "23.66"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "3.34"
$ws.Range("E51").Value = "  +2.64%  "

# Restore the default (General) style so the cells match the original workbook formatting
$dataRange.Style = "Normal"
